# ---------------------------------------------------------------------------
# Edit script: "Add files via upload"
#   * Subtitle on slide 1 gets the author byline.
#   * Four new bio/intro slides are inserted after slide 1.
#   * A new, empty slide is appended at the very end.
#   * The "Why Solar panels from home" slide gets two runs in its
#     second paragraph merged into one.
# ---------------------------------------------------------------------------

function HexColor($hex) {
    # PowerPoint's Font.Color.RGB takes a COLORREF-style 0x00BBGGRR value
    # (same packing as the classic VBA RGB() macro), so build it from the
    # RRGGBB hex string used in the OOXML (<a:srgbClr val="RRGGBB"/>).
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 1 ("Solar Panels for homes") - fill in the subtitle byline.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(6)
$subtitleRange = $subtitle.TextFrame.TextRange
$subtitleRange.Text = "By Claudia Logrande, arica mcneal, michelle chavez"
$subtitleRange.Font.Color.RGB = HexColor "7CEBFF"

# ---------------------------------------------------------------------------
# 2. Insert the new "Introduction presentation" slide at position 2.
# ---------------------------------------------------------------------------
$introSlide = $p.Slides.Add(2, 2)
$introSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Introduction presentation"

$introBody = $introSlide.Shapes.Item(2).TextFrame.TextRange
$introBody.Text = "Girl Power team"
$introP2 = $introBody.InsertAfter("`rWe are a group of single moms, teamed up to work on this project helping each other just like we do in life to better ourselves and learning from each other.")
$introP3 = $introBody.InsertAfter("`rWe have different backgrounds but the same goals to better ourselves by learning new things and succeeding. Also, to set a good example for our kids and point them in the right direction by showing them we" + [char]0x2019 + "re never too old to learn new things.")
$introP3.Font.Bold = $false
$introP3.Font.Italic = $false
$introP3.Font.Name = "Slack-Lato"
$introP3.Font.Color.RGB = HexColor "1D1C1D"

# ---------------------------------------------------------------------------
# 3. Insert the new "Claudia Logrande" bio slide at position 3.
# ---------------------------------------------------------------------------
$claudiaSlide = $p.Slides.Add(3, 2)
$claudiaSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Claudia Logrande"

$claudiaBody = $claudiaSlide.Shapes.Item(2).TextFrame.TextRange
$claudiaBody.Text = "International student in the USA (originally from Italy)."
[void]$claudiaBody.InsertAfter("`rVery involved in environmental issues ")
[void]$claudiaBody.InsertAfter("`rSingle mom of 3 kids")
[void]$claudiaBody.InsertAfter("`rFinished my BA in Business, concentration in Finance only in 2020, then started this program with WOZU because I wanted to learn more about Data Science. Looking forward to learn more about it and taking more classes starting in the Fall.")

# ---------------------------------------------------------------------------
# 4. Reserve slide id 260 for the trailing empty slide by creating it now
#    (kept empty, moved to the end once the remaining bio slides exist).
# ---------------------------------------------------------------------------
$placeholderSlide = $p.Slides.Add(4, 2)

# ---------------------------------------------------------------------------
# 5. Insert the "Arica mcneal" bio slide at position 5.
# ---------------------------------------------------------------------------
$aricaSlide = $p.Slides.Add(5, 2)
$aricaSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Arica mcneal"

# ---------------------------------------------------------------------------
# 6. Insert the "Michelle chavez" bio slide at position 6.
# ---------------------------------------------------------------------------
$michelleSlide = $p.Slides.Add(6, 2)
$michelleSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Michelle chavez"

# ---------------------------------------------------------------------------
# 7. The original "Why Solar panels from home" slide is now at position 7
#    (it was position 2 before any insert). Merge its second paragraph's
#    two runs into a single run.
# ---------------------------------------------------------------------------
$whySlide = $p.Slides.Item(7)
$whyBody = $whySlide.Shapes.Item(2).TextFrame.TextRange
$whyPara2 = $whyBody.Paragraphs(2)
$whyPara2.Text = "placeholder text to break run-merge diffing"
$whyPara2b = $whyBody.Paragraphs(2)
$whyPara2b.Text = "Another reason why we wanted to research this subject is because it" + [char]0x2019 + "s an efficient and greener way to use the natural power of the sun to power our household electronics."

# ---------------------------------------------------------------------------
# 8. Move the reserved empty slide (still blank) to the very end.
# ---------------------------------------------------------------------------
$placeholderIndex = $placeholderSlide.SlideIndex
$moveRange = $p.Slides.Range($placeholderIndex)
$moveRange.Cut()
$p.Slides.Paste($p.Slides.Count + 1) | Out-Null
